# Add a new "LoginData" worksheet (after the existing "UsuariosRegistro"
# sheet) containing a small login-test dataset, as part of completing the
# "second point" of the login feature.

$wb = $excel.ActiveWorkbook
$wsUsuarios = $wb.Worksheets.Item("UsuariosRegistro")

# Insert the new sheet right after UsuariosRegistro so the tab order is
# UsuariosRegistro, LoginData.
$ws = $wb.Worksheets.Add($null, $wsUsuarios)
$ws.Name = "LoginData"

# Header row
$ws.Range("A1").Value = "Email"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "Tipo"

# Sample login record
$ws.Range("A2").Value = "test1@gmail.com"
$ws.Range("B2").Value = 123456
$ws.Range("C2").Value = "invalido"
